$d = $word.ActiveDocument
$d.Content.Find.Execute("politization", $true, $false, $false, $false, $false, $true, 1, $false, "polarization", 2)
